$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# "Handed back: in sync with en-US" -> "Ready for handoff" (Status columns)
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Refreshed handoff-generation timestamps
$overview.Range("G2").Value = "2016-09-06 21:20:41"
$dede.Range("H2").Value = "2016-09-06 21:20:41"
$zhcn.Range("H2").Value = "2016-09-06 21:20:36"

# The Status column no longer needs to be as wide now that the text is
# shorter, so narrow the columns that display it.
$overview.Range("E1").ColumnWidth = 16.333333
$overview.Range("F1").ColumnWidth = 16.333333
$zhcn.Range("C1").ColumnWidth = 16.333333
$dede.Range("C1").ColumnWidth = 16.333333
